$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 29, shifting existing data (rows 29-61) down to rows 31-63.
$ws.Rows.Item(29).Insert()
$ws.Rows.Item(29).Insert()

# Fill row 29 with the new data record.
$ws.Cells.Item(29, 1).Value = 6
$ws.Cells.Item(29, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(29, 3).Value = "Metropolitana"
$ws.Cells.Item(29, 4).Value = 45049
$ws.Cells.Item(29, 5).Value = 13
$ws.Cells.Item(29, 6).Value = "Fruta"
$ws.Cells.Item(29, 7).Value = 100101
$ws.Cells.Item(29, 8).Value = "Berries"
$ws.Cells.Item(29, 9).Value = 100101006
$ws.Cells.Item(29, 10).Value = "Higo"
$ws.Cells.Item(29, 11).Value = "Sin especificar"
$ws.Cells.Item(29, 12).Value = "Primera"
$ws.Cells.Item(29, 13).Value = 80
$ws.Cells.Item(29, 14).Value = 20000
$ws.Cells.Item(29, 15).Value = 20000
$ws.Cells.Item(29, 16).Value = 20000
$ws.Cells.Item(29, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(29, 18).Value = "Región Metropolitana"
$ws.Cells.Item(29, 19).Value = 2857
$ws.Cells.Item(29, 20).Value = 7

# Fill row 30 with the new data record.
$ws.Cells.Item(30, 1).Value = 6
$ws.Cells.Item(30, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(30, 3).Value = "Metropolitana"
$ws.Cells.Item(30, 4).Value = 45049
$ws.Cells.Item(30, 5).Value = 13
$ws.Cells.Item(30, 6).Value = "Fruta"
$ws.Cells.Item(30, 7).Value = 100101
$ws.Cells.Item(30, 8).Value = "Berries"
$ws.Cells.Item(30, 9).Value = 100101006
$ws.Cells.Item(30, 10).Value = "Higo"
$ws.Cells.Item(30, 11).Value = "Sin especificar"
$ws.Cells.Item(30, 12).Value = "Segunda"
$ws.Cells.Item(30, 13).Value = 50
$ws.Cells.Item(30, 14).Value = 15000
$ws.Cells.Item(30, 15).Value = 15000
$ws.Cells.Item(30, 16).Value = 15000
$ws.Cells.Item(30, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(30, 18).Value = "Región Metropolitana"
$ws.Cells.Item(30, 19).Value = 2143
$ws.Cells.Item(30, 20).Value = 7

Write-Host "done"
